$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.07535355310267777
$ws.Range("C2").Value = -0.1246464468973251
$ws.Range("D2").Value = 1.575353553102678
$ws.Range("E2").Value = 1.075353553102692
$ws.Range("F2").Value = -2.224646446897305
$ws.Range("G2").Value = 0.575353553102692
$ws.Range("H2").Value = 0.3753535531026892
$ws.Range("B3").Value = -0.1689259161161161
$ws.Range("C3").Value = 1.531074083883887
$ws.Range("D3").Value = 1.031074083883901
$ws.Range("E3").Value = -2.268925916116096
$ws.Range("F3").Value = 0.531074083883901
$ws.Range("G3").Value = 0.3310740838838982
$ws.Range("B4").Value = 1.844494952278143
$ws.Range("C4").Value = 1.344494952278157
$ws.Range("D4").Value = -1.95550504772184
$ws.Range("E4").Value = 0.8444949522781571
$ws.Range("F4").Value = 0.6444949522781542
$ws.Range("G4").Value = 0.9444949522781514
$ws.Range("H4").Value = 0.6444949522781542
$ws.Range("I4").Value = 0.7444949522781628
$ws.Range("J4").Value = 0.7444949522781628
$ws.Range("B5").Value = 0.746543788682611
$ws.Range("C5").Value = -2.553456211317386
$ws.Range("D5").Value = 0.246543788682611
$ws.Range("E5").Value = 0.04654378868260822
$ws.Range("F5").Value = 0.3465437886826053
$ws.Range("G5").Value = 0.04654378868260822
$ws.Range("H5").Value = 0.1465437886826167
$ws.Range("I5").Value = 0.1465437886826167
$ws.Range("B6").Value = -2.65781059791856
$ws.Range("C6").Value = 0.1421894020814378
$ws.Range("D6").Value = -0.05781059791856499
$ws.Range("E6").Value = 0.2421894020814321
$ws.Range("F6").Value = -0.05781059791856499
$ws.Range("G6").Value = 0.04218940208144348
$ws.Range("H6").Value = 0.04218940208144348
$ws.Range("B7").Value = -0.059948979259083
$ws.Range("C7").Value = -0.2599489792590858
$ws.Range("D7").Value = 0.04005102074091132
$ws.Range("E7").Value = -0.2599489792590858
$ws.Range("F7").Value = -0.1599489792590773
$ws.Range("G7").Value = -0.1599489792590773
$ws.Range("B8").Value = -0.1665024842129383
$ws.Range("C8").Value = 0.1334975157870588
$ws.Range("D8").Value = -0.1665024842129383
$ws.Range("E8").Value = -0.06650248421292981
$ws.Range("F8").Value = -0.06650248421292981
$ws.Range("G8").Value = 0.0334975157870645
$ws.Range("H8").Value = -0.4665024842129213
$ws.Range("I8").Value = -0.3665024842129554
$ws.Range("B9").Value = -0.306554028195724
$ws.Range("C9").Value = -0.6065540281957211
$ws.Range("D9").Value = -0.5065540281957126
$ws.Range("E9").Value = -0.5065540281957126
$ws.Range("F9").Value = -0.4065540281957183
$ws.Range("G9").Value = -0.9065540281957041
$ws.Range("H9").Value = -0.8065540281957382
$ws.Range("B10").Value = -0.176570966028531
$ws.Range("C10").Value = -0.07657096602852248
$ws.Range("D10").Value = -0.07657096602852248
$ws.Range("E10").Value = 0.02342903397147183
$ws.Range("F10").Value = -0.476570966028514
$ws.Range("G10").Value = -0.3765709660285481
$ws.Range("B11").Value = -0.01952633654686198
$ws.Range("C11").Value = -0.01952633654686198
$ws.Range("D11").Value = 0.08047366345313234
$ws.Range("E11").Value = -0.4195263365468535
$ws.Range("F11").Value = -0.3195263365468876
$ws.Range("B12").Value = 0.02083567343742339
$ws.Range("C12").Value = 0.1208356734374177
$ws.Range("D12").Value = -0.3791643265625681
$ws.Range("E12").Value = -0.2791643265626022
$ws.Range("B13").Value = 0.1708795756445881
$ws.Range("C13").Value = -0.3291204243553977
$ws.Range("D13").Value = -0.2291204243554318
$ws.Range("B14").Value = -0.01566635435963382
$ws.Range("C14").Value = 0.08433364564033208
$ws.Range("B15").Value = 0.0449618533261173
